# Update NATMI LR-pair output (Bmp4-Bmpr1b) with re-run TPM results.
# Replaces the existing 6 data rows (2 sending clusters x 3 target clusters)
# with the full 3x3 Cartesian product of sending/target clusters (9 rows),
# reflecting the new TPM-derived statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster (A) / Target cluster (D) labels for the 3x3 grid, plus all
# numeric columns E..T recomputed from the new TPM values.
$rows = @(
    @{ Row=2;  A="ECs";   D="ECs";   E=3; F=1; G=7.620274999999999;  H=22.860825;          I=0.6584612850834004;  J=0.6584612850834003;  K=1; L=0.3333333333333333; M=0.09074700000000001; N=0.272241;           O=0.04501122713837679; P=0.04501122713837678; Q=0.6915170954250001;  R=6.223653858825;     S=0.02963815046471641; T=0.0296381504647164 },
    @{ Row=3;  A="ECs";   D="FAPs";  E=3; F=1; G=7.620274999999999;  H=22.860825;          I=0.6584612850834004;  J=0.6584612850834003;  K=3; L=1;                  M=1.394935666666666;   N=4.184806999999999;  O=0.6918990835593063;  P=0.6918990835593063;  Q=10.62979338730833;   R=95.66814048577497;  S=0.4555887597084878;  T=0.4555887597084878 },
    @{ Row=4;  A="ECs";   D="MuSCs"; E=3; F=1; G=7.620274999999999;  H=22.860825;          I=0.6584612850834004;  J=0.6584612850834003;  K=3; L=1;                  M=0.5304143333333333;  N=1.591243;           O=0.263089689302317;   P=0.263089689302317;   Q=4.041903083941667;   R=36.37712775547499;  S=0.1732343749101962;  T=0.1732343749101962 },
    @{ Row=5;  A="FAPs";  D="ECs";   E=3; F=1; G=3.035834666666667;  H=9.107504;           I=0.262323813236933;   J=0.262323813236933;   K=1; L=0.3333333333333333; M=0.09074700000000001; N=0.272241;           O=0.04501122713837679; P=0.04501122713837678; Q=0.275492888496;      R=2.479435996464;     S=0.01180751674141273; T=0.01180751674141272 },
    @{ Row=6;  A="FAPs";  D="FAPs";  E=3; F=1; G=3.035834666666667;  H=9.107504;           I=0.262323813236933;   J=0.262323813236933;   K=3; L=1;                  M=1.394935666666666;   N=4.184806999999999;  O=0.6918990835593063;  P=0.6918990835593063;  Q=4.234794054636444;   R=38.113146491728;    S=0.1815016059744166;  T=0.1815016059744166 },
    @{ Row=7;  A="FAPs";  D="MuSCs"; E=3; F=1; G=3.035834666666667;  H=9.107504;           I=0.262323813236933;   J=0.262323813236933;   K=3; L=1;                  M=0.5304143333333333;  N=1.591243;           O=0.263089689302317;   P=0.263089689302317;   Q=1.610250220830222;   R=14.492251987472;    S=0.06901469052110375; T=0.06901469052110375 },
    @{ Row=8;  A="MuSCs"; D="ECs";   E=3; F=1; G=0.9167423333333334; H=2.750227;           I=0.07921490167966665; J=0.07921490167966663; K=1; L=0.3333333333333333; M=0.09074700000000001; N=0.272241;           O=0.04501122713837679; P=0.04501122713837678; Q=0.08319161652300001; R=0.7487245487070001; S=0.00356555993224766; T=0.00356555993224766 },
    @{ Row=9;  A="MuSCs"; D="FAPs";  E=3; F=1; G=0.9167423333333334; H=2.750227;           I=0.07921490167966665; J=0.07921490167966663; K=3; L=1;                  M=1.394935666666666;   N=4.184806999999999;  O=0.6918990835593063;  P=0.6918990835593063;  Q=1.278796577909889;   R=11.509169201189;    S=0.0548087178764019;  T=0.0548087178764019 },
    @{ Row=10; A="MuSCs"; D="MuSCs"; E=3; F=1; G=0.9167423333333334; H=2.750227;           I=0.07921490167966665; J=0.07921490167966663; K=3; L=1;                  M=0.5304143333333333;  N=1.591243;           O=0.263089689302317;   P=0.263089689302317;   Q=0.4862532735734445;  R=4.376279462161;     S=0.02084062387101709; T=0.02084062387101708 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A        # A: Sending cluster
    $ws.Cells.Item($row, 2).Value = "Bmp4"      # B: Ligand symbol (unchanged)
    $ws.Cells.Item($row, 3).Value = "Bmpr1b"    # C: Receptor symbol (unchanged)
    $ws.Cells.Item($row, 4).Value = $r.D        # D: Target cluster
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
